$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data collected 3 Sep 2025 (serial date 45903) for plots IP3, C3, IP2, C2, C1, IP1
$plots = @("IP3", "C3", "IP2", "C2", "C1", "IP1")
$dates = @(45903, 45903, 45903, 45903, 45903, 45903)
$times = @(0.50486111111111109, 0.50972222222222219, 0.52777777777777779, 0.53055555555555556, 0.53333333333333333, 0.53611111111111109)
$headspace1 = @(31.4, 31.2, 32.700000000000003, 31.5, 31.5, 31.5)
$headspace2 = @(3, 3, 3, 3, 3, 3)
$extLen = @(2, 2, 2, 2, 2, 2)
$duration = @(3, 4, 3, 3, 3, 3)

for ($i = 0; $i -lt 6; $i++) {
    $row = 38 + $i

    $ws.Cells.Item($row, 1).Value = $plots[$i]
    $ws.Cells.Item($row, 2).Value = $dates[$i]
    $ws.Cells.Item($row, 3).Value = $times[$i]
    $ws.Cells.Item($row, 4).Value = $headspace1[$i]
    $ws.Cells.Item($row, 5).Value = $headspace2[$i]
    $ws.Cells.Item($row, 9).Value = $extLen[$i]
    $ws.Cells.Item($row, 10).Value = $duration[$i]
}

# Update the saved view/selection to reflect the newly added rows
$window = $excel.ActiveWindow
$window.ScrollRow = 28
$window.ScrollColumn = 5
$ws.Range("J44").Select()

$wb.Save()
